# Make sheet2 ("Sheet2 - Numbers") match reader/sheet2: extend the data
# with a new AA column (values 100..129 next to the existing D/K columns),
# which pushes the sheet's used range out to D1:AA30, and leave the sheet
# selected/active with AA1:AA30 highlighted (AA1 as the active cell) -
# mirroring how Excel recorded the edit session.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2 - Numbers")

# New column AA, rows 1-30: 100, 101, 102, ... 129
for ($i = 0; $i -lt 30; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 27).Value = 100 + $i
}

# Sheet2 becomes the active/selected sheet (was Sheet4 before the edit),
# with the new AA column selected - AA1 is the active cell, AA1:AA30 is
# the highlighted range.
$ws.Select() | Out-Null
$ws.Range("AA1:AA30").Select() | Out-Null

# Scroll the window so column AA is close to the left edge of the view
# (best effort - mirrors the recorded topLeftCell="O1" viewport).
$excel.ActiveWindow.ScrollColumn = 15
$excel.ActiveWindow.ScrollRow = 1

# Sheet4's page setup was recorded with an explicit paper size (A4) on
# the re-save.
$ws4 = $wb.Worksheets.Item("Sheet4 - Dates")
$ws4.PageSetup.PaperSize = 9
